$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 922.44446
$ws.Range("I6").Value = 922.44446
$ws.Range("K6").Value = 2767.33338
$ws.Range("M6").Value = -2655.33338
$ws.Range("H17").Value = 554.6053000000001
$ws.Range("J17").Value = 302.5
$ws.Range("L17").Value = 907.5
$ws.Range("N17").Value = -1243.5
$ws.Range("H64").Value = 3309.1538
$ws.Range("I64").Value = 3032.2222
$ws.Range("J64").Value = 3932.25
$ws.Range("K64").Value = 3032.2222
$ws.Range("L64").Value = 3932.25
$ws.Range("M64").Value = -2784.2222
$ws.Range("N64").Value = -4428.25
$ws.Range("H67").Value = 3309.1538
$ws.Range("I67").Value = 3032.2222
$ws.Range("J67").Value = 3932.25
$ws.Range("K67").Value = 3032.2222
$ws.Range("L67").Value = 3932.25
$ws.Range("M67").Value = -2174.2222
$ws.Range("N67").Value = -5648.25
$ws.Range("H93").Value = 37530.81
$ws.Range("J93").Value = 37530.81
$ws.Range("L93").Value = 37530.81
$ws.Range("N93").Value = -42522.81
$ws.Range("H98").Value = 5293.467
$ws.Range("I98").Value = 3279.15
$ws.Range("J98").Value = 6904.92
$ws.Range("K98").Value = 3279.15
$ws.Range("L98").Value = 6904.92
$ws.Range("M98").Value = -1781.15
$ws.Range("N98").Value = -9900.92
$ws.Range("H116").Value = 338262.66
$ws.Range("I116").Value = 835356.4399999999
$ws.Range("J116").Value = 6866.8335
$ws.Range("K116").Value = 835356.4399999999
$ws.Range("L116").Value = 6866.8335
$ws.Range("M116").Value = -831914.4399999999
$ws.Range("N116").Value = -13750.8335
$ws.Range("H122").Value = 5293.467
$ws.Range("I122").Value = 3279.15
$ws.Range("J122").Value = 6904.92
$ws.Range("K122").Value = 9837.450000000001
$ws.Range("L122").Value = 20714.76
$ws.Range("M122").Value = -7387.450000000001
$ws.Range("N122").Value = -25614.76
$ws.Range("H123").Value = 39419.09
$ws.Range("J123").Value = 39419.09
$ws.Range("L123").Value = 39419.09
$ws.Range("N123").Value = -49219.09
$ws.Range("H132").Value = 40563296
$ws.Range("I132").Value = 45639124
$ws.Range("J132").Value = 3340555.2
$ws.Range("K132").Value = 136917372
$ws.Range("L132").Value = 10021665.6
$ws.Range("M132").Value = -136914842
$ws.Range("N132").Value = -10026725.6
$ws.Range("H137").Value = 664267.25
$ws.Range("I137").Value = 1539127.4
$ws.Range("J137").Value = 2787.7073
$ws.Range("K137").Value = 4617382.199999999
$ws.Range("L137").Value = 8363.1219
$ws.Range("M137").Value = -4614832.199999999
$ws.Range("N137").Value = -13463.1219
$ws.Range("H138").Value = 2949.541
$ws.Range("I138").Value = 1583.7222
$ws.Range("J138").Value = 3521.279
$ws.Range("K138").Value = 4751.1666
$ws.Range("L138").Value = 10563.837
$ws.Range("M138").Value = 388.8334000000004
$ws.Range("N138").Value = -20843.837

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3586.7722
$ws.Range("I32").Value = 3235.9253
$ws.Range("J32").Value = 5545.6665
$ws.Range("K32").Value = 3235.9253
$ws.Range("L32").Value = 5545.6665
$ws.Range("M32").Value = -2948.9253
$ws.Range("N32").Value = -6119.6665
$ws.Range("H103").Value = 34857.145
$ws.Range("J103").Value = 34857.145
$ws.Range("L103").Value = 34857.145
$ws.Range("N103").Value = -37201.145
$ws.Range("H132").Value = 5457
$ws.Range("I132").Value = 3433.1667
$ws.Range("J132").Value = 9504.666999999999
$ws.Range("K132").Value = 10299.5001
$ws.Range("L132").Value = 28514.001
$ws.Range("M132").Value = -7769.500100000001
$ws.Range("N132").Value = -33574.001
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140
$ws.Range("H137").Value = 39886
$ws.Range("J137").Value = 39886
$ws.Range("L137").Value = 39886
$ws.Range("N137").Value = -50086

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2407.8
$ws.Range("I99").Value = 1554.625
$ws.Range("J99").Value = 3382.8572
$ws.Range("K99").Value = 1554.625
$ws.Range("L99").Value = 3382.8572
$ws.Range("M99").Value = -56.625
$ws.Range("N99").Value = -6378.8572
$ws.Range("H103").Value = 34819.23
$ws.Range("J103").Value = 34819.23
$ws.Range("L103").Value = 34819.23
$ws.Range("N103").Value = -37163.23
$ws.Range("H132").Value = 53488.723
$ws.Range("J132").Value = 53488.723
$ws.Range("L132").Value = 53488.723
$ws.Range("N132").Value = -63608.723
$ws.Range("H134").Value = 3748.804
$ws.Range("I134").Value = 1293.909
$ws.Range("J134").Value = 5611.1377
$ws.Range("K134").Value = 3881.727
$ws.Range("L134").Value = 16833.4131
$ws.Range("M134").Value = -1346.727
$ws.Range("N134").Value = -21903.4131
$ws.Range("H135").Value = 48756.617
$ws.Range("J135").Value = 48756.617
$ws.Range("L135").Value = 48756.617
$ws.Range("N135").Value = -58896.617

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4403.3794
$ws.Range("I132").Value = 3273.9473
$ws.Range("K132").Value = 9821.841899999999
$ws.Range("M132").Value = -7291.841899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H68").Value = 2254.6707
$ws.Range("I68").Value = 770.3929000000001
$ws.Range("J68").Value = 2983.7896
$ws.Range("K68").Value = 2311.1787
$ws.Range("L68").Value = 8951.3688
$ws.Range("M68").Value = -1500.1787
$ws.Range("N68").Value = -10573.3688
$ws.Range("H71").Value = 2254.6707
$ws.Range("I71").Value = 770.3929000000001
$ws.Range("J71").Value = 2983.7896
$ws.Range("K71").Value = 6933.5361
$ws.Range("L71").Value = 26854.1064
$ws.Range("M71").Value = -2877.5361
$ws.Range("N71").Value = -34966.1064
$ws.Range("H107").Value = 15317.871
$ws.Range("J107").Value = 32032.121
$ws.Range("L107").Value = 96096.363
$ws.Range("N107").Value = -99936.363
$ws.Range("H113").Value = 5000824.5
$ws.Range("I113").Value = 671.0769
$ws.Range("J113").Value = 10417657
$ws.Range("K113").Value = 2013.2307
$ws.Range("L113").Value = 31252971
$ws.Range("M113").Value = 156.7692999999999
$ws.Range("N113").Value = -31257311
$ws.Range("H122").Value = 2216.4827
$ws.Range("I122").Value = 758.7778
$ws.Range("J122").Value = 2872.45
$ws.Range("K122").Value = 6829.000199999999
$ws.Range("L122").Value = 25852.05
$ws.Range("M122").Value = -4379.000199999999
$ws.Range("N122").Value = -30752.05

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 6004
$ws.Range("I9").Value = 4604
$ws.Range("J9").Value = 9504
$ws.Range("K9").Value = 4604
$ws.Range("L9").Value = 9504
$ws.Range("M9").Value = -4434
$ws.Range("N9").Value = -9844
$ws.Range("H102").Value = 4365.067
$ws.Range("I102").Value = 3500
$ws.Range("J102").Value = 5662.6665
$ws.Range("K102").Value = 3500
$ws.Range("L102").Value = 5662.6665
$ws.Range("M102").Value = -1878
$ws.Range("N102").Value = -8906.666499999999
$ws.Range("H111").Value = 35293
$ws.Range("J111").Value = 35293
$ws.Range("L111").Value = 35293
$ws.Range("N111").Value = -41427
$ws.Range("H122").Value = 5725.9
$ws.Range("I122").Value = 4666.3335
$ws.Range("J122").Value = 6180
$ws.Range("K122").Value = 13999.0005
$ws.Range("L122").Value = 18540
$ws.Range("M122").Value = -11549.0005
$ws.Range("N122").Value = -23440
$ws.Range("H124").Value = 41865.715
$ws.Range("J124").Value = 41865.715
$ws.Range("L124").Value = 41865.715
$ws.Range("N124").Value = -51685.715
$ws.Range("H135").Value = 48334.285
$ws.Range("J135").Value = 48334.285
$ws.Range("L135").Value = 48334.285
$ws.Range("N135").Value = -58474.285

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H93").Value = 6946990
$ws.Range("I93").Value = 8549295
$ws.Range("J93").Value = 3668
$ws.Range("K93").Value = 8549295
$ws.Range("L93").Value = 3668
$ws.Range("M93").Value = -8548047
$ws.Range("N93").Value = -6164
$ws.Range("H115").Value = 30500
$ws.Range("J115").Value = 30500
$ws.Range("L115").Value = 30500
$ws.Range("N115").Value = -32850
$ws.Range("H127").Value = 24450.938
$ws.Range("J127").Value = 24450.938
$ws.Range("L127").Value = 24450.938
$ws.Range("N127").Value = -34370.93799999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 40048.668
$ws.Range("J39").Value = 40048.668
$ws.Range("L39").Value = 40048.668
$ws.Range("N39").Value = -40874.668
$ws.Range("H49").Value = 20016000
$ws.Range("I49").Value = 50005000
$ws.Range("J49").Value = 23333
$ws.Range("K49").Value = 50005000
$ws.Range("L49").Value = 23333
$ws.Range("M49").Value = -50004770
$ws.Range("N49").Value = -23793
$ws.Range("H118").Value = 29223.334
$ws.Range("J118").Value = 29223.334
$ws.Range("L118").Value = 29223.334
$ws.Range("N118").Value = -32537.334
$ws.Range("H123").Value = 35008.566
$ws.Range("J123").Value = 35008.566
$ws.Range("L123").Value = 35008.566
$ws.Range("N123").Value = -44808.566
$ws.Range("H126").Value = 535284.3
$ws.Range("I126").Value = 1913.7142
$ws.Range("J126").Value = 822483.9
$ws.Range("K126").Value = 5741.142599999999
$ws.Range("L126").Value = 2467451.7
$ws.Range("M126").Value = -3271.142599999999
$ws.Range("N126").Value = -2472391.7
$ws.Range("H132").Value = 2950.0557
$ws.Range("I132").Value = 1767
$ws.Range("J132").Value = 4133.1113
$ws.Range("K132").Value = 5301
$ws.Range("L132").Value = 12399.3339
$ws.Range("M132").Value = -2771
$ws.Range("N132").Value = -17459.3339
$ws.Range("H136").Value = 5002.15
$ws.Range("I136").Value = 2424.25
$ws.Range("J136").Value = 6720.75
$ws.Range("K136").Value = 7272.75
$ws.Range("L136").Value = 20162.25
$ws.Range("M136").Value = -4722.75
$ws.Range("N136").Value = -25262.25
$ws.Range("H137").Value = 35960
$ws.Range("J137").Value = 35960
$ws.Range("L137").Value = 35960
$ws.Range("N137").Value = -46160
